$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = ";djsavlkcsadvbsdavgbshdklvbfdsjklvhbsdjkal"
$ws.Range("A3").Value = "cdscdwvwrvfreavcsfdvfsav"

$ws.Range("A4").Select()
